$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price values that would otherwise be mis-parsed as numbers by Excel
# (single-decimal-point strings) - force text format before assigning, so exact
# string formatting (trailing zeros, leading zeros, full precision) is preserved.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.42"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.80"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.22"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.29"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.08"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.92"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "495.12"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000151"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.47"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.35"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.19"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.94"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.00"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.87"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "442.41"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.85"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.88"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.47"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.38"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.61"

# D-column price values with two decimal separators - Excel cannot parse these as
# numbers, so they stay text automatically; no NumberFormat override needed.
$ws.Range("D2").Value = "68.659.72"
$ws.Range("D3").Value = "3.751.22"
$ws.Range("D7").Value = "3.749.94"
$ws.Range("D15").Value = "4.377.40"
$ws.Range("D16").Value = "3.748.67"
$ws.Range("D17").Value = "68.694.02"
$ws.Range("D34").Value = "3.897.08"
$ws.Range("D35").Value = "3.685.75"
$ws.Range("D49").Value = "2.820.30"

# E-column volume/percentage change values (always contain "%", stay text automatically).
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("E21").Value = "  +20.38%  "
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("E24").Value = "  +7.66%  "
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("E30").Value = "  +7.51%  "
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("E32").Value = "  +3.46%  "
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("E51").Value = "  +3.12%  "
